$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set NumberFormat to Text for the Price column range so numeric-looking
# strings (e.g. "1.002") are preserved as text instead of being coerced
# into numbers, then restore the default "Normal" style so no stray
# style index is left on the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.010.36"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.861.36"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "312.22"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.5117"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("D8").Value = "0.3850"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("D9").Value = "0.08278"
$ws.Range("E9").Value = "  -8.29%  "
$ws.Range("D10").Value = "1.113"
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").Value = "41.57"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("E12").Value = "  -2.51%  "
$ws.Range("D13").Value = "20.59"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").Value = "1.864.05"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").Value = "7.253"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").Value = "90.69"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "0.06657"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").Value = "17.68"
$ws.Range("E20").Value = "  -3.02%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "6.009"
$ws.Range("E22").Value = "  -1.93%  "
$ws.Range("D23").Value = "28.032.11"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "11.08"
$ws.Range("E24").Value = "  -3.10%  "
$ws.Range("D25").Value = "2.230"
$ws.Range("E25").Value = "  -1.79%  "
$ws.Range("D26").Value = "2.073.19"
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("D27").Value = "2.515"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("D28").Value = "157.57"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").Value = "20.49"
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("D30").Value = "124.82"
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("E32").Value = "  -3.28%  "
$ws.Range("D33").Value = "5.947"
$ws.Range("E33").Value = "  +5.74%  "
$ws.Range("D34").Value = "3.596"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Value = "9.382"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").Value = "0.02411"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "0.06491"
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("D38").Value = "0.2177"
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("D39").Value = "0.6603"
$ws.Range("E39").Value = "  +2.93%  "
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("D41").Value = "5.013"
$ws.Range("E41").Value = "  +1.82%  "
$ws.Range("D42").Value = "1.227"
$ws.Range("E42").Value = "  -4.38%  "
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("D44").Value = "0.6155"
$ws.Range("E44").Value = "  +1.69%  "
$ws.Range("D45").Value = "13.02"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("D46").Value = "1.280"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").Value = "3.657"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").Value = "2.011"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("D49").Value = "1.208"
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("D50").Value = "119.93"
$ws.Range("E50").Value = "  -1.23%  "
$ws.Range("D51").Value = "78.94"
$ws.Range("E51").Value = "  -1.15%  "

# Restore default style on the price column (removes the temporary text
# number-format so the cells match the original unstyled appearance).
$priceRange.Style = "Normal"
